$wb = $excel.ActiveWorkbook

# --- "SensorBox consumption" sheet: no cell/data changes, just loses tab focus ---
# --- "Energy saving example" sheet: update inputs and become the active sheet ---
$wsEnergy = $wb.Worksheets.Item("Energy saving example")

# Update input cells (dependent formulas recalc automatically)
$wsEnergy.Range("I3").Value = 250
$wsEnergy.Range("C6").Value = 12
$wsEnergy.Range("F17").Value = 1.5

# Make this sheet the active one (moves tabSelected from "SensorBox consumption"
# to "Energy saving example" and updates workbook.xml's activeTab)
$wsEnergy.Activate()

# Update the selected cell on the now-active sheet
$wsEnergy.Range("C7").Select()
